$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.178.10"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.357.96"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "543.16"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "133.56"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  +4.94%  "
$ws.Range("D9").Value = "0.106"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "23.82"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.775.42"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "58.122.08"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "2.365.47"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "330.25"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "6.83"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").Value = "63.53"
$ws.Range("E23").Value = "  +2.67%  "
$ws.Range("D24").Value = "0.166"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "8.23"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("D27").Value = "1.32"
$ws.Range("E27").Value = "  -5.65%  "
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "170.97"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").Value = "0.417"
$ws.Range("E39").Value = "  +9.99%  "
$ws.Range("D40").Value = "142.10"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("D41").Value = "3.66"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").Value = "288.38"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "0.0949"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").Value = "0.0516"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "18.95"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.566"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "0.0222"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").Value = "0.384"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Value = "4.71"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("E51").Value = "  +0.55%  "
